# Atualização de bases das ligas, do dia: 10-06-2024 às 07:08
# The odds-feed refresh re-keyed several fixtures: rows 181/182 swap their
# match data (id stays put in column A), and rows 313/314/315 rotate their
# match data in a 3-cycle (313<-315, 314<-313, 315<-314).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $Values) {
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

# --- Row 181 <-- old Row 182 data (id/A column stays 179) ---
Set-Row 181 @{
    B  = 6810169
    E  = "Westerlo"
    F  = "Cercle Brugge"
    G  = 4
    H  = 2
    I  = 4
    J  = 1
    K  = "H"
    L  = 3.4
    M  = 3.8
    N  = 1.95
    O  = 3.3
    P  = 3.6
    Q  = 2.05
    R  = 0.25
    S  = 2
    T  = 1.85
    U  = 2.75
    V  = 1.975
    W  = 1.875
    X  = 2.3
    Y  = -1
    AA = 1
    AB = -1
    AC = 0.9750000000000001
}

# --- Row 182 <-- old Row 181 data (id/A column stays 180) ---
Set-Row 182 @{
    B  = 6810167
    E  = "Club Brugge"
    F  = "KV Kortrijk"
    G  = 3
    H  = 3
    I  = 1
    J  = 0
    K  = "D"
    L  = 1.125
    M  = 8.5
    N  = 17
    O  = 1.125
    P  = 8.5
    Q  = 17
    R  = -2.25
    S  = 1.85
    T  = 2
    U  = 3.5
    V  = 1.9
    W  = 1.95
    X  = -1
    Y  = 7.5
    AA = -1
    AB = 1
    AC = 0.8999999999999999
}

# --- Row 313 <-- old Row 315 data (id/A column stays 311) ---
Set-Row 313 @{
    B  = 8009904
    E  = "Antwerp"
    F  = "Anderlecht"
    G  = 3
    H  = 1
    L  = 3.5
    M  = 3.5
    N  = 1.833
    O  = 3.25
    P  = 3.5
    Q  = 1.9
    R  = 0.5
    S  = 1.85
    T  = 2
    U  = 3
    V  = 2.025
    W  = 1.825
    X  = 2.25
    AA = 0.8500000000000001
    AC = 1.025
    AD = -1
}

# --- Row 314 <-- old Row 313 data (id/A column stays 312) ---
Set-Row 314 @{
    B  = 8009865
    E  = "Union Saint Gilloise"
    F  = "Genk"
    G  = 2
    K  = "H"
    L  = 1.666
    M  = 3.75
    N  = 4.333
    O  = 1.5
    P  = 4.2
    Q  = 5
    R  = -1
    S  = 1.875
    T  = 1.975
    V  = 1.875
    W  = 1.975
    X  = 0.5
    Y  = -1
    AA = 0.875
    AB = -1
    AD = 0.9750000000000001
}

# --- Row 315 <-- old Row 314 data (id/A column stays 313) ---
Set-Row 315 @{
    B  = 8009325
    E  = "Club Brugge"
    F  = "Cercle Brugge"
    G  = 0
    H  = 0
    K  = "D"
    L  = 1.444
    M  = 4.5
    N  = 5.5
    O  = 1.615
    P  = 3.8
    Q  = 4.5
    R  = -0.75
    S  = 1.825
    T  = 2.025
    U  = 3.25
    V  = 2.05
    W  = 1.8
    X  = -1
    Y  = 2.8
    AA = -1
    AB = 1.025
    AC = -1
    AD = 0.8
}
